$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix misspelled / inconsistent names ---
# "Alison Bell" -> "Alisson Bell" (matches the spelling already used elsewhere, e.g. A162)
$ws.Range("B55").Value = "Alisson Bell"

# "Sasha Cheplin" -> "Sasha Chepelin" (matches the spelling already used elsewhere, e.g. A194)
$ws.Range("B64").Value = "Sasha Chepelin"

# "Louis McMillan" -> "Louis MacMillan" (matches the spelling already used elsewhere, e.g. B90)
$ws.Range("A151").Value = "Louis MacMillan"

# "Gay Barber" -> "Dan Barber" (matches the name already used elsewhere, e.g. A104 / B146)
$ws.Range("A161").Value = "Dan Barber"

# --- Add two new rows of data at the bottom of the table ---
# Row 219: copy formatting from row 217 (name / name / pulled / hh / hh)
$ws.Range("A217:E217").Copy()
$ws.Range("A219:E219").PasteSpecial(-4122)
$ws.Range("A219").Value = "Dan O'Riordan"
$ws.Range("B219").Value = "Maja Thomson"
$ws.Range("C219").Value = "pulled"
$ws.Range("D219").Value = "hh"
$ws.Range("E219").Value = "hh"

# Row 220: copy formatting from row 218 (name / name / slept_with / hh / hh)
$ws.Range("A218:E218").Copy()
$ws.Range("A220:E220").PasteSpecial(-4122)
$ws.Range("A220").Value = "Euan Patton"
$ws.Range("B220").Value = "Maja Thomson"
$ws.Range("C220").Value = "slept_with"
$ws.Range("D220").Value = "hh"
$ws.Range("E220").Value = "hh"

[void]($excel.CutCopyMode = $false)

# --- Restore the frozen header row / original selection (unchanged by this edit) ---
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
[void]$ws.Range("B3").Select()

Write-Output "edits applied"
